$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "10÷5="
$t.Cell(1,2).Range.Text = "40÷9="
$t.Cell(1,3).Range.Text = "53÷4="
$t.Cell(1,4).Range.Text = "91÷7="
$t.Cell(1,5).Range.Text = "84÷6="
$t.Cell(5,1).Range.Text = "23÷9="
$t.Cell(5,2).Range.Text = "73÷4="
$t.Cell(5,3).Range.Text = "92÷7="
$t.Cell(5,4).Range.Text = "27÷9="
$t.Cell(5,5).Range.Text = "54÷5="
$t.Cell(9,1).Range.Text = "91÷3="
$t.Cell(9,2).Range.Text = "99÷8="
$t.Cell(9,3).Range.Text = "46÷2="
$t.Cell(9,4).Range.Text = "92÷9="
$t.Cell(9,5).Range.Text = "90÷3="
$t.Cell(13,1).Range.Text = "13÷9="
$t.Cell(13,2).Range.Text = "81÷4="
$t.Cell(13,3).Range.Text = "11÷4="
$t.Cell(13,4).Range.Text = "65÷3="
$t.Cell(13,5).Range.Text = "45÷6="
$t.Cell(17,1).Range.Text = "52÷5="
$t.Cell(17,2).Range.Text = "54÷8="
$t.Cell(17,3).Range.Text = "56÷5="
$t.Cell(17,4).Range.Text = "51÷9="
$t.Cell(17,5).Range.Text = "21÷5="
